# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 70
$ws.Range("H70").Value = 3321.2222
$ws.Range("I70").Value = 3798
$ws.Range("J70").Value = 2939.8
$ws.Range("K70").Value = 11394
$ws.Range("L70").Value = 8819.400000000001
$ws.Range("M70").Value = -11124
$ws.Range("N70").Value = -9359.400000000001
# ALC row 73
$ws.Range("H73").Value = 3321.2222
$ws.Range("I73").Value = 3798
$ws.Range("J73").Value = 2939.8
$ws.Range("K73").Value = 11394
$ws.Range("L73").Value = 8819.400000000001
$ws.Range("M73").Value = -10458
$ws.Range("N73").Value = -10691.4
# ALC row 80
$ws.Range("H80").Value = 608.08
$ws.Range("I80").Value = 544.2857
$ws.Range("J80").Value = 689.2727
$ws.Range("K80").Value = 1632.8571
$ws.Range("L80").Value = 2067.8181
$ws.Range("M80").Value = -634.8571000000002
$ws.Range("N80").Value = -4063.8181
# ALC row 83
$ws.Range("H83").Value = 608.08
$ws.Range("I83").Value = 544.2857
$ws.Range("J83").Value = 689.2727
$ws.Range("K83").Value = 4898.571300000001
$ws.Range("L83").Value = 6203.454299999999
$ws.Range("M83").Value = 93.42869999999948
$ws.Range("N83").Value = -16187.4543
# ALC row 96
$ws.Range("H96").Value = 464.66666
$ws.Range("I96").Value = 447
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 1341
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = 32
$ws.Range("N96").Value = -4246

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = -474
# ARM row 32
$ws.Range("H32").Value = 4156.3647
$ws.Range("I32").Value = 2966.0657
$ws.Range("J32").Value = 9741.615
$ws.Range("K32").Value = 2966.0657
$ws.Range("L32").Value = 9741.615
$ws.Range("M32").Value = -2679.0657
$ws.Range("N32").Value = -10315.615
# ARM row 61
$ws.Range("H61").Value = 5759.5
$ws.Range("I61").Value = 6128.1304
$ws.Range("K61").Value = 6128.1304
$ws.Range("M61").Value = -5916.1304
# ARM row 133
$ws.Range("H133").Value = 34925
$ws.Range("J133").Value = 34925
$ws.Range("L133").Value = 34925
$ws.Range("N133").Value = -39985
# ARM row 136
$ws.Range("H136").Value = 5759.5
$ws.Range("I136").Value = 6128.1304
$ws.Range("K136").Value = 18384.3912
$ws.Range("M136").Value = -15834.3912

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -480
# BSM row 94
$ws.Range("H94").Value = 2312
$ws.Range("I94").Value = 1893.7142
$ws.Range("J94").Value = 2800
$ws.Range("K94").Value = 1893.7142
$ws.Range("L94").Value = 2800
$ws.Range("M94").Value = -1442.7142
$ws.Range("N94").Value = -3702
# BSM row 123
$ws.Range("H123").Value = 39796
$ws.Range("J123").Value = 39796
$ws.Range("L123").Value = 39796
$ws.Range("N123").Value = -49596
# BSM row 125
$ws.Range("H125").Value = 52779.332
$ws.Range("J125").Value = 52779.332
$ws.Range("L125").Value = 52779.332
$ws.Range("N125").Value = -62619.332
# BSM row 127
$ws.Range("H127").Value = 56390
$ws.Range("J127").Value = 56390
$ws.Range("L127").Value = 56390
$ws.Range("N127").Value = -66310
# BSM row 129
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

$ws = $wb.Worksheets.Item("CRP")
# CRP row 6
$ws.Range("H6").Value = 36666868
$ws.Range("I6").Value = 55000000
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 55000000
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -54999887
$ws.Range("N6").Value = -826
# CRP row 52
$ws.Range("H52").Value = 28310
$ws.Range("J52").Value = 28310
$ws.Range("L52").Value = 28310
$ws.Range("N52").Value = -28898
# CRP row 132
$ws.Range("H132").Value = 3759.8462
$ws.Range("I132").Value = 3542.2222
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 10626.6666
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -8096.6666
$ws.Range("N132").Value = -17808.5
# CRP row 137
$ws.Range("H137").Value = 35290
$ws.Range("J137").Value = 35290
$ws.Range("L137").Value = 35290
$ws.Range("N137").Value = -45490

$ws = $wb.Worksheets.Item("CUL")
# CUL row 37
$ws.Range("H37").Value = 40470.59
$ws.Range("J37").Value = 40470.59
$ws.Range("L37").Value = 121411.77
$ws.Range("N37").Value = -121635.77
# CUL row 68
$ws.Range("H68").Value = 3996.6667
$ws.Range("J68").Value = 5750
$ws.Range("L68").Value = 17250
$ws.Range("N68").Value = -18872
# CUL row 71
$ws.Range("H71").Value = 3996.6667
$ws.Range("J71").Value = 5750
$ws.Range("L71").Value = 51750
$ws.Range("N71").Value = -59862
# CUL row 97
$ws.Range("H97").Value = 16667200
$ws.Range("I97").Value = 25000326
$ws.Range("J97").Value = 950
$ws.Range("K97").Value = 75000978
$ws.Range("L97").Value = 2850
$ws.Range("M97").Value = -75000482
$ws.Range("N97").Value = -3842
# CUL row 131
$ws.Range("H131").Value = 2174838
$ws.Range("I131").Value = 10000500
$ws.Range("J131").Value = 1043.0834
$ws.Range("K131").Value = 30001500
$ws.Range("L131").Value = 3129.2502
$ws.Range("M131").Value = -29996460
$ws.Range("N131").Value = -13209.2502

$ws = $wb.Worksheets.Item("GSM")
# GSM row 114
$ws.Range("H114").Value = 49831
$ws.Range("J114").Value = 49831
$ws.Range("L114").Value = 49831
$ws.Range("N114").Value = -58509
# GSM row 137
$ws.Range("H137").Value = 29786.666
$ws.Range("J137").Value = 29786.666
$ws.Range("L137").Value = 29786.666
$ws.Range("N137").Value = -39986.666

$ws = $wb.Worksheets.Item("LTW")
# LTW row 68
$ws.Range("H68").Value = 90910800
$ws.Range("I68").Value = 1642
$ws.Range("J68").Value = 166668430
$ws.Range("K68").Value = 1642
$ws.Range("L68").Value = 166668430
$ws.Range("M68").Value = -893
$ws.Range("N68").Value = -166669928
# LTW row 71
$ws.Range("H71").Value = 90910800
$ws.Range("I71").Value = 1642
$ws.Range("J71").Value = 166668430
$ws.Range("K71").Value = 8210
$ws.Range("L71").Value = 833342150
$ws.Range("M71").Value = -4466
$ws.Range("N71").Value = -833349638
# LTW row 132
$ws.Range("H132").Value = 12386493
$ws.Range("I132").Value = 14948209
$ws.Range("K132").Value = 44844627
$ws.Range("M132").Value = -44842097
# LTW row 134
$ws.Range("H134").Value = 35000
$ws.Range("J134").Value = 35000
$ws.Range("L134").Value = 35000
$ws.Range("N134").Value = -45140
# LTW row 136
$ws.Range("H136").Value = 6145.857
$ws.Range("I136").Value = 6878.2173
$ws.Range("K136").Value = 20634.6519
$ws.Range("M136").Value = -18084.6519

$ws = $wb.Worksheets.Item("WVR")
# WVR row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
# WVR row 132
$ws.Range("H132").Value = 1490.591
$ws.Range("I132").Value = 988.6667
$ws.Range("K132").Value = 2966.0001
$ws.Range("M132").Value = -436.0001000000002
# WVR row 136
$ws.Range("H136").Value = 1364.0476
$ws.Range("I136").Value = 861.6667
$ws.Range("J136").Value = 2033.8889
$ws.Range("K136").Value = 2585.0001
$ws.Range("L136").Value = 6101.6667
$ws.Range("M136").Value = -35.0001000000002
$ws.Range("N136").Value = -11201.6667

